$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 203.55556
$ws.Range("I38").Value = 203.55556
$ws.Range("K38").Value = 610.66668
$ws.Range("M38").Value = -238.66668
$ws.Range("H39").Value = 519.7692
$ws.Range("I39").Value = 141.54546
$ws.Range("K39").Value = 424.63638
$ws.Range("M39").Value = -128.63638
$ws.Range("H64").Value = 6674.4736
$ws.Range("I64").Value = 3870.7778
$ws.Range("J64").Value = 9197.799999999999
$ws.Range("K64").Value = 3870.7778
$ws.Range("L64").Value = 9197.799999999999
$ws.Range("M64").Value = -3622.7778
$ws.Range("N64").Value = -9693.799999999999
$ws.Range("H67").Value = 6674.4736
$ws.Range("I67").Value = 3870.7778
$ws.Range("J67").Value = 9197.799999999999
$ws.Range("K67").Value = 3870.7778
$ws.Range("L67").Value = 9197.799999999999
$ws.Range("M67").Value = -3012.7778
$ws.Range("N67").Value = -10913.8
$ws.Range("H74").Value = 7566.269
$ws.Range("I74").Value = 5680.3335
$ws.Range("J74").Value = 8564.706
$ws.Range("K74").Value = 5680.3335
$ws.Range("L74").Value = 8564.706
$ws.Range("M74").Value = -4744.3335
$ws.Range("N74").Value = -10436.706
$ws.Range("H77").Value = 7566.269
$ws.Range("I77").Value = 5680.3335
$ws.Range("J77").Value = 8564.706
$ws.Range("K77").Value = 28401.6675
$ws.Range("L77").Value = 42823.53
$ws.Range("M77").Value = -23721.6675
$ws.Range("N77").Value = -52183.53
$ws.Range("H100").Value = 3980.3928
$ws.Range("I100").Value = 1632.3684
$ws.Range("K100").Value = 1632.3684
$ws.Range("M100").Value = -1091.3684
$ws.Range("H137").Value = 2562.8635
$ws.Range("I137").Value = 1210
$ws.Range("J137").Value = 2776.4736
$ws.Range("K137").Value = 3630
$ws.Range("L137").Value = 8329.4208
$ws.Range("M137").Value = -1080
$ws.Range("N137").Value = -13429.4208
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1805.3429
$ws.Range("I2").Value = 1681.9354
$ws.Range("J2").Value = 2761.75
$ws.Range("K2").Value = 1681.9354
$ws.Range("L2").Value = 2761.75
$ws.Range("M2").Value = -1568.9354
$ws.Range("N2").Value = -2987.75
$ws.Range("H63").Value = 11250
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 11250
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H116").Value = 1805.3429
$ws.Range("I116").Value = 1681.9354
$ws.Range("J116").Value = 2761.75
$ws.Range("K116").Value = 1681.9354
$ws.Range("L116").Value = 2761.75
$ws.Range("M116").Value = 612.0645999999999
$ws.Range("N116").Value = -7349.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1805.3429
$ws.Range("I3").Value = 1681.9354
$ws.Range("J3").Value = 2761.75
$ws.Range("K3").Value = 1681.9354
$ws.Range("L3").Value = 2761.75
$ws.Range("M3").Value = -1567.9354
$ws.Range("N3").Value = -2989.75
$ws.Range("H107").Value = 2245.606
$ws.Range("I107").Value = 2340.3667
$ws.Range("K107").Value = 2340.3667
$ws.Range("M107").Value = -420.3667
$ws.Range("H133").Value = 78232.8
$ws.Range("J133").Value = 70055
$ws.Range("L133").Value = 70055
$ws.Range("N133").Value = -80175
$ws.Range("H134").Value = 693522.5600000001
$ws.Range("I134").Value = 803971.5600000001
$ws.Range("K134").Value = 2411914.68
$ws.Range("M134").Value = -2409379.68
$ws.Range("H139").Value = 47499.5
$ws.Range("I139").Value = 39999
$ws.Range("J139").Value = 55000
$ws.Range("K139").Value = 39999
$ws.Range("L139").Value = 55000
$ws.Range("M139").Value = -34859
$ws.Range("N139").Value = -65280
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 825750.25
$ws.Range("I58").Value = 1236797.9
$ws.Range("K58").Value = 1236797.9
$ws.Range("M58").Value = -1236594.9
$ws.Range("H99").Value = 4870.643
$ws.Range("I99").Value = 5410
$ws.Range("J99").Value = 3899.8
$ws.Range("K99").Value = 5410
$ws.Range("L99").Value = 3899.8
$ws.Range("M99").Value = -3912
$ws.Range("N99").Value = -6895.8
$ws.Range("H122").Value = 972
$ws.Range("I122").Value = 964.6667
$ws.Range("K122").Value = 2894.0001
$ws.Range("M122").Value = -444.0001000000002
$ws.Range("H126").Value = 4870.643
$ws.Range("I126").Value = 5410
$ws.Range("J126").Value = 3899.8
$ws.Range("K126").Value = 16230
$ws.Range("L126").Value = 11699.4
$ws.Range("M126").Value = -13760
$ws.Range("N126").Value = -16639.4
$ws.Range("H132").Value = 9274346
$ws.Range("I132").Value = 18625.904
$ws.Range("K132").Value = 55877.712
$ws.Range("M132").Value = -53347.712
$ws.Range("H136").Value = 825750.25
$ws.Range("I136").Value = 1236797.9
$ws.Range("K136").Value = 3710393.7
$ws.Range("M136").Value = -3707843.7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 95546060
$ws.Range("I4").Value = 185000060
$ws.Range("J4").Value = 21001064
$ws.Range("K4").Value = 555000180
$ws.Range("L4").Value = 63003192
$ws.Range("M4").Value = -555000068
$ws.Range("N4").Value = -63003416
$ws.Range("H40").Value = 119
$ws.Range("I40").Value = 126.5
$ws.Range("J40").Value = 111.5
$ws.Range("K40").Value = 506
$ws.Range("L40").Value = 446
$ws.Range("M40").Value = -437
$ws.Range("N40").Value = -584
$ws.Range("H107").Value = 546.4545000000001
$ws.Range("J107").Value = 601.7143
$ws.Range("L107").Value = 1805.1429
$ws.Range("N107").Value = -5645.1429
$ws.Range("H113").Value = 1703.8148
$ws.Range("I113").Value = 1156.125
$ws.Range("J113").Value = 1934.421
$ws.Range("K113").Value = 3468.375
$ws.Range("L113").Value = 5803.263
$ws.Range("M113").Value = -1298.375
$ws.Range("N113").Value = -10143.263
$ws.Range("H122").Value = 741.8
$ws.Range("I122").Value = 638.9
$ws.Range("J122").Value = 947.6
$ws.Range("K122").Value = 5750.099999999999
$ws.Range("L122").Value = 8528.4
$ws.Range("M122").Value = -3300.099999999999
$ws.Range("N122").Value = -13428.4
$ws.Range("H129").Value = 2606.375
$ws.Range("I129").Value = 1397.1428
$ws.Range("J129").Value = 3546.889
$ws.Range("K129").Value = 4191.428400000001
$ws.Range("L129").Value = 10640.667
$ws.Range("M129").Value = 808.5715999999993
$ws.Range("N129").Value = -20640.667
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2951.4075
$ws.Range("I113").Value = 1609.7693
$ws.Range("K113").Value = 1609.7693
$ws.Range("M113").Value = 560.2307000000001
$ws.Range("H132").Value = 2040.1333
$ws.Range("I132").Value = 1917
$ws.Range("K132").Value = 5751
$ws.Range("M132").Value = -3221
$ws.Range("H141").Value = 69420
$ws.Range("J141").Value = 69420
$ws.Range("L141").Value = 69420
$ws.Range("N141").Value = -79780
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2259
$ws.Range("I100").Value = 2307
$ws.Range("J100").Value = 2244
$ws.Range("K100").Value = 4614
$ws.Range("L100").Value = 4488
$ws.Range("M100").Value = -4073
$ws.Range("N100").Value = -5570
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H107").Value = 979.3953
$ws.Range("I107").Value = 708.43335
$ws.Range("K107").Value = 2125.30005
$ws.Range("M107").Value = -205.3000499999998
$ws.Range("H113").Value = 1329.8334
$ws.Range("I113").Value = 1021.8182
$ws.Range("J113").Value = 1813.8572
$ws.Range("K113").Value = 3065.4546
$ws.Range("L113").Value = 5441.571599999999
$ws.Range("M113").Value = -895.4546
$ws.Range("N113").Value = -9781.571599999999
$ws.Range("H122").Value = 2038.3611
$ws.Range("I122").Value = 1550.2963
$ws.Range("J122").Value = 3502.5557
$ws.Range("K122").Value = 4650.8889
$ws.Range("L122").Value = 10507.6671
$ws.Range("M122").Value = -2200.8889
$ws.Range("N122").Value = -15407.6671
$ws.Range("H126").Value = 2397.111
$ws.Range("I126").Value = 2505.1365
$ws.Range("J126").Value = 2227.3572
$ws.Range("K126").Value = 7515.4095
$ws.Range("L126").Value = 6682.071599999999
$ws.Range("M126").Value = -5045.4095
$ws.Range("N126").Value = -11622.0716
$ws.Range("H140").Value = 49214.5
$ws.Range("J140").Value = 49214.5
$ws.Range("L140").Value = 49214.5
$ws.Range("N140").Value = -59574.5
$ws.Range("H141").Value = 90000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 90000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 90000
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -100360
